# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner market-price refresh values to the Leviathan_Profits workbook
# (static cached values only -- no formulas in these sheets)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 490
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 750
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 750
$ws.Range("M55").Value = 114
$ws.Range("N55").Value = -1178

$ws.Range("H58").Value = 777.7273
$ws.Range("I58").Value = 438.125
$ws.Range("J58").Value = 1683.3334
$ws.Range("K58").Value = 1314.375
$ws.Range("L58").Value = 5050.0002
$ws.Range("M58").Value = -1164.375
$ws.Range("N58").Value = -5350.0002

$ws.Range("H61").Value = 1283.5714
$ws.Range("I61").Value = 1393.6
$ws.Range("K61").Value = 4180.799999999999
$ws.Range("M61").Value = -4008.799999999999

$ws.Range("H99").Value = 250291120
$ws.Range("I99").Value = 82255.5
$ws.Range("K99").Value = 246766.5
$ws.Range("M99").Value = -245268.5

$ws.Range("H113").Value = 4680
$ws.Range("I113").Value = 3266.6667
$ws.Range("K113").Value = 3266.6667
$ws.Range("M113").Value = -12.66670000000022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41795.805
$ws.Range("I32").Value = 25716.781
$ws.Range("K32").Value = 25716.781
$ws.Range("M32").Value = -25429.781

$ws.Range("H45").Value = 16687.867
$ws.Range("I45").Value = 15000.125
$ws.Range("J45").Value = 18616.715
$ws.Range("K45").Value = 15000.125
$ws.Range("L45").Value = 18616.715
$ws.Range("M45").Value = -14623.125
$ws.Range("N45").Value = -19370.715

$ws.Range("H110").Value = 900.25
$ws.Range("I110").Value = 868.6667
$ws.Range("K110").Value = 868.6667
$ws.Range("M110").Value = 1176.3333

$ws.Range("H122").Value = 25486
$ws.Range("I122").Value = 25486
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 76458
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -74008
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6178099
$ws.Range("I20").Value = 10104903
$ws.Range("K20").Value = 10104903
$ws.Range("M20").Value = -10104656

$ws.Range("H97").Value = 13085.333
$ws.Range("I97").Value = 13085.333
$ws.Range("K97").Value = 13085.333
$ws.Range("M97").Value = -12094.333

$ws.Range("H99").Value = 1251.909
$ws.Range("I99").Value = 1251.909
$ws.Range("K99").Value = 1251.909
$ws.Range("M99").Value = 246.0909999999999

$ws.Range("H102").Value = 23537.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1589.4375
$ws.Range("I58").Value = 1379.3077
$ws.Range("K58").Value = 1379.3077
$ws.Range("M58").Value = -1176.3077

$ws.Range("H93").Value = 25241.857
$ws.Range("I93").Value = 23159.4
$ws.Range("K93").Value = 23159.4
$ws.Range("M93").Value = -21287.4

$ws.Range("H103").Value = 142878940
$ws.Range("I103").Value = 166687090
$ws.Range("K103").Value = 166687090
$ws.Range("M103").Value = -166685918

$ws.Range("H132").Value = 2161.7334
$ws.Range("I132").Value = 2173.4644
$ws.Range("K132").Value = 6520.3932
$ws.Range("M132").Value = -3990.3932

$ws.Range("H134").Value = 2253.0667
$ws.Range("I134").Value = 2230.1365
$ws.Range("K134").Value = 6690.4095
$ws.Range("M134").Value = -4155.4095

$ws.Range("H136").Value = 1589.4375
$ws.Range("I136").Value = 1379.3077
$ws.Range("K136").Value = 4137.9231
$ws.Range("M136").Value = -1587.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 6264.6665
$ws.Range("I63").Value = 1897.5
$ws.Range("K63").Value = 5692.5
$ws.Range("M63").Value = -4943.5

$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 6264.6665
$ws.Range("I66").Value = 1897.5
$ws.Range("K66").Value = 17077.5
$ws.Range("M66").Value = -13333.5

$ws.Range("H94").Value = 11795.143
$ws.Range("I94").Value = 5036.25
$ws.Range("J94").Value = 14498.7
$ws.Range("K94").Value = 15108.75
$ws.Range("L94").Value = 43496.10000000001
$ws.Range("M94").Value = -14432.75
$ws.Range("N94").Value = -44848.10000000001

$ws.Range("H97").Value = 1854.875
$ws.Range("I97").Value = 1366.3334
$ws.Range("J97").Value = 2148
$ws.Range("K97").Value = 4099.0002
$ws.Range("L97").Value = 6444
$ws.Range("M97").Value = -3603.0002
$ws.Range("N97").Value = -7436

$ws.Range("H98").Value = 608.8333
$ws.Range("I98").Value = 491
$ws.Range("J98").Value = 667.75
$ws.Range("K98").Value = 1473
$ws.Range("L98").Value = 2003.25
$ws.Range("M98").Value = 25
$ws.Range("N98").Value = -4999.25

$ws.Range("H117").Value = 2232.75
$ws.Range("J117").Value = 3232.6
$ws.Range("L117").Value = 9697.799999999999
$ws.Range("N117").Value = -16581.8

$ws.Range("H129").Value = 135949.67
$ws.Range("J129").Value = 3128.8333
$ws.Range("L129").Value = 9386.499899999999
$ws.Range("N129").Value = -19386.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H80").Value = 7565.7393
$ws.Range("I80").Value = 12794.3
$ws.Range("K80").Value = 12794.3
$ws.Range("M80").Value = -11796.3

$ws.Range("H83").Value = 7565.7393
$ws.Range("I83").Value = 12794.3
$ws.Range("K83").Value = 63971.5
$ws.Range("M83").Value = -58979.5

$ws.Range("H99").Value = 12394.143
$ws.Range("I99").Value = 9373.333000000001
$ws.Range("K99").Value = 9373.333000000001
$ws.Range("M99").Value = -7127.333000000001

$ws.Range("H102").Value = 2200
$ws.Range("I102").Value = 1750
$ws.Range("K102").Value = 1750
$ws.Range("M102").Value = -128

$ws.Range("H113").Value = 55558996
$ws.Range("I113").Value = 83335990
$ws.Range("J113").Value = 5007
$ws.Range("K113").Value = 83335990
$ws.Range("L113").Value = 5007
$ws.Range("M113").Value = -83333820
$ws.Range("N113").Value = -9347

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1711.84
$ws.Range("I22").Value = 1586.7333
$ws.Range("J22").Value = 1899.5
$ws.Range("K22").Value = 1586.7333
$ws.Range("L22").Value = 1899.5
$ws.Range("M22").Value = -1291.7333
$ws.Range("N22").Value = -2489.5

$ws.Range("H27").Value = 1711.84
$ws.Range("I27").Value = 1586.7333
$ws.Range("J27").Value = 1899.5
$ws.Range("K27").Value = 1586.7333
$ws.Range("L27").Value = 1899.5
$ws.Range("M27").Value = -1479.7333
$ws.Range("N27").Value = -2113.5

$ws.Range("H55").Value = 414.57895
$ws.Range("I55").Value = 312.1111
$ws.Range("J55").Value = 506.8
$ws.Range("K55").Value = 312.1111
$ws.Range("L55").Value = 506.8
$ws.Range("M55").Value = -139.1111
$ws.Range("N55").Value = -852.8

$ws.Range("H61").Value = 16181.625
$ws.Range("I61").Value = 20258.834
$ws.Range("K61").Value = 20258.834
$ws.Range("M61").Value = -20056.834

$ws.Range("H68").Value = 2467
$ws.Range("I68").Value = 2467
$ws.Range("K68").Value = 2467
$ws.Range("M68").Value = -1718

$ws.Range("H71").Value = 2467
$ws.Range("I71").Value = 2467
$ws.Range("K71").Value = 12335
$ws.Range("M71").Value = -8591

$ws.Range("H113").Value = 16181.625
$ws.Range("I113").Value = 20258.834
$ws.Range("K113").Value = 20258.834
$ws.Range("M113").Value = -18088.834

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 132212.38
$ws.Range("I62").Value = 8116.6665
$ws.Range("K62").Value = 8116.6665
$ws.Range("M62").Value = -7492.6665

$ws.Range("H65").Value = 132212.38
$ws.Range("I65").Value = 8116.6665
$ws.Range("K65").Value = 40583.3325
$ws.Range("M65").Value = -37463.3325

$ws.Range("H107").Value = 29413006
$ws.Range("I107").Value = 1291.5
$ws.Range("K107").Value = 3874.5
$ws.Range("M107").Value = -1954.5

$ws.Range("H113").Value = 1439.6
